# "Reorganización completa: limpieza de módulos antiguos, nuevas entregas y optimización"
#
# plantilla_animales.xlsx:
#  - rename sheet "Datos" -> "animales"
#  - drop the old 24-column (A:X) header row (with its bold/blue/centered
#    header style) in favor of a leaner 19-column (A:S) header row using
#    plain snake_case names and no special formatting
#  - a handful of old columns (Lote, Sector, Grupo, Vendedor, Procedencia)
#    are retired in the process

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "animales"

# New (shorter) header list for A1:S1
$headers = @(
    "codigo",
    "nombre",
    "tipo_ingreso",
    "sexo",
    "fecha_nacimiento",
    "fecha_compra",
    "finca",
    "raza",
    "potrero",
    "peso_nacimiento",
    "peso_compra",
    "precio_compra",
    "salud",
    "estado",
    "inventariado",
    "color",
    "hierro",
    "condicion_corporal",
    "comentario"
)

# The old sheet had 24 columns (A:X); the new layout only needs 19 (A:S),
# so remove the trailing five (old Lote/Sector/Grupo/Vendedor/Procedencia
# area is gone / folded away).
$ws.Columns("T:X").Delete()

# Overwrite the header row with the new plain-text labels and strip the old
# bold/white-on-blue centered header styling back to the default style.
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.ClearFormats()
}
